$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in D3:D18, preserving cell formatting/style
$ws.Range("D3:D18").ClearContents()

# Update the selection to match the target state
$ws.Range("D3:D18").Select()
